# Updates cryptos list (prices / 1h volume %, plus a few re-ranked coins)
# with fresh GitHub Actions scrape data.
#
# Note: several "Price" values look numeric (e.g. "303.38") but must stay
# text, matching the sheet's existing inlineStr/string cells. A leading
# apostrophe forces Excel to store them as text (quote-prefixed) instead of
# silently parsing them into Double values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.202.84'
$ws.Range("E2").Value = '  +1.82%  '
$ws.Range("D3").Value = '2.383.02'
$ws.Range("E3").Value = '  +4.13%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''303.38'
$ws.Range("E5").Value = '  +0.98%  '
$ws.Range("D6").Value = '''97.63'
$ws.Range("E6").Value = '  +3.05%  '
$ws.Range("D7").Value = '''0.509'
$ws.Range("E7").Value = '  +0.69%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '''0.502'
$ws.Range("E9").Value = '  +2.43%  '
$ws.Range("D10").Value = '''34.34'
$ws.Range("E10").Value = '  +0.01%  '
$ws.Range("E11").Value = '  +1.31%  '
$ws.Range("E12").Value = '  +2.45%  '
$ws.Range("E13").Value = '  -2.59%  '
$ws.Range("D14").Value = '''6.82'
$ws.Range("E14").Value = '  +2.06%  '
$ws.Range("D15").Value = '2.752.97'
$ws.Range("E15").Value = '  +4.01%  '
$ws.Range("D16").Value = '2.403.43'
$ws.Range("E16").Value = '  +5.12%  '
$ws.Range("D17").Value = '''0.809'
$ws.Range("E17").Value = '  +4.30%  '
$ws.Range("D18").Value = '43.180.50'
$ws.Range("E18").Value = '  +1.89%  '
$ws.Range("D19").Value = '''12.22'
$ws.Range("E19").Value = '  +0.70%  '
$ws.Range("E20").Value = '  +5.82%  '
$ws.Range("E21").Value = '  +0.46%  '
$ws.Range("D22").Value = '''68.49'
$ws.Range("E22").Value = '  +1.53%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").Value = '''235.58'
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("B24").Value = 'ImmutableX'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D24").Value = '''2.25'
$ws.Range("E24").Value = '  +0.54%  '
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").Value = '''2.45'
$ws.Range("E25").Value = '  +2.30%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("E27").Value = '  +3.38%  '
$ws.Range("E28").Value = '  +0.36%  '
$ws.Range("D29").Value = '''9.16'
$ws.Range("E29").Value = '  +1.80%  '
$ws.Range("D30").Value = '''31.77'
$ws.Range("E30").Value = '  +0.74%  '
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("E32").Value = '  +3.05%  '
$ws.Range("D33").Value = '''0.0741'
$ws.Range("E33").Value = '  +7.31%  '
$ws.Range("D34").Value = '''17.28'
$ws.Range("E34").Value = '  -0.50%  '
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").Value = '''0.105'
$ws.Range("E35").Value = '  +5.36%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '''1.85'
$ws.Range("E36").Value = '  +7.34%  '
$ws.Range("E37").Value = '  -0.94%  '
$ws.Range("E38").Value = '  -1.29%  '
$ws.Range("E39").Value = '  +4.91%  '
$ws.Range("D40").Value = '''22.52'
$ws.Range("E40").Value = '  +14.77%  '
$ws.Range("E41").Value = '  +0.70%  '
$ws.Range("D42").Value = '1.960.27'
$ws.Range("E42").Value = '  +0.83%  '
$ws.Range("D43").Value = '''103.52'
$ws.Range("E43").Value = '  -36.84%  '
$ws.Range("E44").Value = '  +1.48%  '
$ws.Range("E45").Value = '  +2.37%  '
$ws.Range("E46").Value = '  +1.57%  '
$ws.Range("D47").Value = '''9.22'
$ws.Range("E47").Value = '  -10.37%  '
$ws.Range("B48").Value = 'MultiversX'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D48").Value = '''52.86'
$ws.Range("E48").Value = '  +0.42%  '
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").Value = '''1.51'
$ws.Range("E49").Value = '  +3.38%  '
$ws.Range("B50").Value = 'BitcoinSV'
$ws.Range("C50").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D50").Value = '''72.00'
$ws.Range("E50").Value = '  +2.55%  '
$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").Value = '''1.15'
$ws.Range("E51").Value = '  +1.94%  '
